$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New LDA block (rows 18-25), mirroring the GLM Fit / QDA blocks above ---

# Row 18: date + method name "LDA" + first data row (Top 5)
$ws.Range("A18").HorizontalAlignment = -4108
$ws.Range("A18").NumberFormat = "d-mmm-yy"
$ws.Range("A18").Value = 45949
$ws.Range("B18").Value = "LDA"
$ws.Range("C18").Value = "Top 5"
$ws.Range("D18").Value = 0.68179999999999996

# Rows 19-25: remaining Top N rows for LDA
$ws.Range("C19").Value = "Top 10"
$ws.Range("D19").Value = 0.77270000000000005

$ws.Range("C20").Value = "Top 15"
$ws.Range("D20").Value = 0.77270000000000005

$ws.Range("C21").Value = "Top 20"
$ws.Range("D21").Value = 0.77270000000000005

$ws.Range("C22").Value = "Top 25"
$ws.Range("D22").Value = 0.81820000000000004

$ws.Range("C23").Value = "Top 30"
$ws.Range("D23").Value = 0.77270000000000005

$ws.Range("C24").Value = "Top 35"
$ws.Range("D24").Value = 0.77270000000000005

$ws.Range("C25").Value = "Top 37"
$ws.Range("D25").Value = 0.72729999999999995

# --- Bold the date cells and method-name cells (A2/A10/A18, B2/B10/B18) ---
$ws.Range("A2").Font.Bold = $true
$ws.Range("A10").Font.Bold = $true
$ws.Range("A18").Font.Bold = $true

$ws.Range("B2").Font.Bold = $true
$ws.Range("B10").Font.Bold = $true
$ws.Range("B18").Font.Bold = $true

# --- Update the selection to match the saved workbook (row 2 selected) ---
$ws.Range("A2:XFD2").Select()
